$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: Invalid / Absent -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count / Real -> 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count / Real -> 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Absent -> 1
$ws.Range("H6").Value = 1

# Row 7: Absent -> 1
$ws.Range("H7").Value = 1

# Row 8: Absent -> 1
$ws.Range("H8").Value = 1

# Row 9: Total Attendance Count / Real -> 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: Absent -> 1
$ws.Range("H10").Value = 1

# Row 11: Absent -> 1
$ws.Range("H11").Value = 1

# Row 12: Absent -> 1
$ws.Range("H12").Value = 1

# Row 13: Total Attendance Count / Real -> 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

# Row 14: Absent -> 1
$ws.Range("H14").Value = 1

# Row 15: Absent -> 1
$ws.Range("H15").Value = 1

# Row 16: Absent -> 1
$ws.Range("H16").Value = 1

# Row 17: Absent -> 1
$ws.Range("H17").Value = 1

# Row 18: Absent -> 1
$ws.Range("H18").Value = 1
